# "redbus count number of buses"
#
# Adds a new "redbus" worksheet (after the existing "irctc" sheet) with a
# from/to/date header row and one data row (bengaluru -> chennai), mirroring
# the layout/style already used on the "irctc" sheet. Also moves the active
# selection: "irctc" keeps a plain A1:C1 range selection (no longer the
# active tab), while the new "redbus" sheet becomes the active tab with
# A4 selected.

$wb = $excel.ActiveWorkbook

$irctc = $wb.Worksheets.Item(1)

# New sheet, inserted right after "irctc" so tab order is irctc, redbus.
$redbus = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $irctc)
$redbus.Name = "redbus"

# Header row - reuses the same shared strings ("from"/"to"/"date") as irctc.
$redbus.Range("A1").Value = "from"
$redbus.Range("B1").Value = "to"
$redbus.Range("C1").Value = "date"

# Highlight the header row the same way irctc's header row is highlighted.
$redbus.Range("A1:C1").Interior.Color = 65535

# Data row. Write "chennai" (B2) before "bengaluru" (A2) so the shared
# string table gains them in that order.
$redbus.Range("B2").Value = "chennai"
$redbus.Range("A2").Value = "bengaluru"

# Column A sized to fit its content.
$redbus.Columns.Item(1).ColumnWidth = 10

# irctc: selection becomes a plain header-row range, no longer the active tab.
$irctc.Range("A1:C1").Select()

# redbus becomes the active sheet/tab, with A4 selected.
$redbus.Activate()
$redbus.Range("A4").Select()
